$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "29.431.58"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.61%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.966.78"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "  +3.55%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.09%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "326.21"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.08%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.08%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.39%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3909"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.23%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 2)
$c.NumberFormat = "@"
$c.Value = "OKB"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "46.19"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.90%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 2)
$c.NumberFormat = "@"
$c.Value = "Dogecoin"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.07918"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.48%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 2)
$c.NumberFormat = "@"
$c.Value = "Polygon"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.9869"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.23%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 2)
$c.NumberFormat = "@"
$c.Value = "Solana"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "22.73"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.12%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 2)
$c.NumberFormat = "@"
$c.Value = "WrappedEther"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "1.992.18"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "  +3.98%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 2)
$c.NumberFormat = "@"
$c.Value = "Chainlink"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "7.167"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.47%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 2)
$c.NumberFormat = "@"
$c.Value = "Polkadot"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "5.810"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.38%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 2)
$c.NumberFormat = "@"
$c.Value = "TRON"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.07081"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.32%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 2)
$c.NumberFormat = "@"
$c.Value = "Litecoin"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "87.64"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.75%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 2)
$c.NumberFormat = "@"
$c.Value = "BinanceUSD"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.13%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 2)
$c.NumberFormat = "@"
$c.Value = "ShibaInu"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.000009906"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.42%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 2)
$c.NumberFormat = "@"
$c.Value = "Avalanche"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "17.24"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.06%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 2)
$c.NumberFormat = "@"
$c.Value = "Dai"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.09%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 2)
$c.NumberFormat = "@"
$c.Value = "WrappedBTC"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "29.429.97"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.59%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 2)
$c.NumberFormat = "@"
$c.Value = "Uniswap"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "5.520"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.21%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 2)
$c.NumberFormat = "@"
$c.Value = "Cosmos"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "11.12"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.40%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 2)
$c.NumberFormat = "@"
$c.Value = "WrappedliquidstakedEther2.0"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.207.38"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.76%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 2)
$c.NumberFormat = "@"
$c.Value = "Toncoin"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "2.102"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.21%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 2)
$c.NumberFormat = "@"
$c.Value = "Monero"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "157.99"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.40%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 2)
$c.NumberFormat = "@"
$c.Value = "EthereumClassic"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "19.46"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.08%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 2)
$c.NumberFormat = "@"
$c.Value = "InternetComputer(DFINITY)"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "5.771"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.04%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 2)
$c.NumberFormat = "@"
$c.Value = "BitcoinCash"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "119.33"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.71%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 2)
$c.NumberFormat = "@"
$c.Value = "LidoDAOToken"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "1.895"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.33%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 2)
$c.NumberFormat = "@"
$c.Value = "Stellar"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.09406"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.47%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 2)
$c.NumberFormat = "@"
$c.Value = "ImmutableX"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.8894"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.44%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 2)
$c.NumberFormat = "@"
$c.Value = "Filecoin"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "5.225"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.59%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 2)
$c.NumberFormat = "@"
$c.Value = "ARBITRUM"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.316"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.52%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 2)
$c.NumberFormat = "@"
$c.Value = "HuobiToken"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "3.165"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.32%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 2)
$c.NumberFormat = "@"
$c.Value = "Hedera"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.05804"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.42%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 2)
$c.NumberFormat = "@"
$c.Value = "TrustWalletToken"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "1.166"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.92%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 2)
$c.NumberFormat = "@"
$c.Value = "VeChain"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.02100"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.69%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 2)
$c.NumberFormat = "@"
$c.Value = "FraxShare"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "7.724"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.13%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 2)
$c.NumberFormat = "@"
$c.Value = "TheSandbox"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.5702"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.08%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 2)
$c.NumberFormat = "@"
$c.Value = "PEPE"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.000003131"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "  +48.06%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 2)
$c.NumberFormat = "@"
$c.Value = "Algorand"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.1792"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.32%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 2)
$c.NumberFormat = "@"
$c.Value = "Aptos"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "9.635"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.77%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 2)
$c.NumberFormat = "@"
$c.Value = "MXToken"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "2.747"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "  +7.05%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 2)
$c.NumberFormat = "@"
$c.Value = "EnergySwap"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "11.76"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.14%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 2)
$c.NumberFormat = "@"
$c.Value = "Decentraland"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.5331"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.46%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 2)
$c.NumberFormat = "@"
$c.Value = "RenderToken"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "2.176"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.39%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 2)
$c.NumberFormat = "@"
$c.Value = "Cronos"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.06910"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.59%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "1.823"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.35%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 2)
$c.NumberFormat = "@"
$c.Value = "Quant"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "113.22"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.21%  "
$c.NumberFormat = "General"
$c.Style = "Normal"
